# Update save-load system and some player interactions
$wb = $excel.ActiveWorkbook

# --- enemies sheet: fix swapped gold-drop values for wild dog / goblin ---
$enemies = $wb.Worksheets.Item("enemies")
$enemies.Range("F4").Value = 7
$enemies.Range("F5").Value = 3

# --- locations sheet: rename stray "monsterlvl" header back to "lvl" ---
$locations = $wb.Worksheets.Item("locations")
$locations.Range("B1").Value = "lvl"

# --- update per-sheet selections left over from editing ---
[void]$locations.Range("B4").Select()

# enemies becomes the active / selected sheet (was "loot" before)
[void]$enemies.Activate()
[void]$enemies.Range("F4").Select()
